$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GOLD")

# Update the last cell's value (was "Automation", now "lastRow")
$ws.Range("E4").Value = "lastRow"

# Move the active selection to E4 to match the saved view state
$ws.Range("E4").Select()
